# Apply the Celerio account.xlsx template update to the "Search" sheet:
#  1. Insert a new search-criteria row ("search_full_text") right after the
#     "search_criteria" row (i.e. before the former row 4), shifting all the
#     following rows down by one.
#  2. Set the new row's label/value cells.
#  3. Replace the former "account_homeAddress"/"Role" rows (now at the
#     bottom) with a single "Security Roles"/"securityRoles" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Search")

# 1. Insert a new row before row 4 (pushes rows 4..13 down to 5..14,
#    including the C/D/E birth-date-range cells on the old row 11).
$ws.Rows("4:4").Insert()

# 2. Populate the newly inserted row 4 with the full-text search criteria.
$ws.Range("A4").Value = '${msg.getProperty(''search_full_text'')}'
$ws.Range("B4").Value = '${search_full_text}'

# 3. The old "account_homeAddress" row is now row 14 and the old "Role" row
#    was pushed out to row 15. Replace row 14 with the new "Security Roles"
#    row and delete the now-superfluous row 15.
$ws.Range("A14").Value = '${msg.getProperty(''Security Roles'')}'
$ws.Range("B14").Value = '${securityRoles}'
$ws.Rows("15:15").Delete()
